$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.694.24"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "2.018.12"
$ws.Range("E3").Value = "  -4.19%  "
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "333.05"
$ws.Range("E5").Value = "  -3.57%  "
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "0.5061"
$ws.Range("E7").Value = "  -3.13%  "
$ws.Range("D8").Value = "0.4269"
$ws.Range("E8").Value = "  -3.84%  "
$ws.Range("D9").Value = "54.22"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").Value = "0.09248"
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("D11").Value = "1.131"
$ws.Range("E11").Value = "  -3.72%  "
$ws.Range("D12").Value = "23.68"
$ws.Range("E12").Value = "  -5.72%  "
$ws.Range("D13").Value = "8.160"
$ws.Range("E13").Value = "  -6.46%  "
$ws.Range("D14").Value = "2.004.61"
$ws.Range("E14").Value = "  -5.14%  "
$ws.Range("D15").Value = "6.576"
$ws.Range("E15").Value = "  -4.99%  "
$ws.Range("D16").Value = "96.14"
$ws.Range("E16").Value = "  -5.53%  "
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("E18").Value = "  -2.95%  "
$ws.Range("D19").Value = "0.06668"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "20.03"
$ws.Range("E20").Value = "  -6.07%  "
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "6.028"
$ws.Range("E22").Value = "  -4.59%  "
$ws.Range("D23").Value = "29.704.57"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").Value = "12.07"
$ws.Range("E24").Value = "  -4.52%  "
$ws.Range("D25").Value = "2.278"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").Value = "159.94"
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("D27").Value = "20.88"
$ws.Range("E27").Value = "  -5.21%  "
$ws.Range("D28").Value = "6.465"
$ws.Range("E28").Value = "  -5.68%  "
$ws.Range("D29").Value = "2.363"
$ws.Range("E29").Value = "  -7.17%  "
$ws.Range("D30").Value = "129.12"
$ws.Range("E30").Value = "  -3.24%  "
$ws.Range("D31").Value = "1.067"
$ws.Range("E31").Value = "  -6.97%  "
$ws.Range("D32").Value = "1.600"
$ws.Range("E32").Value = "  -9.36%  "
$ws.Range("D33").Value = "0.09987"
$ws.Range("E33").Value = "  -5.39%  "
$ws.Range("D34").Value = "5.907"
$ws.Range("E34").Value = "  -5.63%  "
$ws.Range("D35").Value = "3.803"
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("D36").Value = "9.665"
$ws.Range("E36").Value = "  -8.35%  "
$ws.Range("D37").Value = "0.02482"
$ws.Range("E37").Value = "  -5.70%  "
$ws.Range("D38").Value = "1.326"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").Value = "0.06422"
$ws.Range("E39").Value = "  -5.54%  "
$ws.Range("D40").Value = "0.6617"
$ws.Range("E40").Value = "  -6.06%  "
$ws.Range("D41").Value = "11.88"
$ws.Range("E41").Value = "  -5.40%  "
$ws.Range("D42").Value = "0.2094"
$ws.Range("E42").Value = "  -6.23%  "
$ws.Range("D43").Value = "1.009"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").Value = "0.6401"
$ws.Range("E44").Value = "  -6.35%  "
$ws.Range("D45").Value = "13.66"
$ws.Range("E45").Value = "  -5.80%  "
$ws.Range("D46").Value = "2.227"
$ws.Range("E46").Value = "  -5.76%  "
$ws.Range("E47").Value = "  -5.17%  "
$ws.Range("D48").Value = "3.540"
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("D49").Value = "0.07031"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("D50").Value = "1.144"
$ws.Range("E50").Value = "  -4.89%  "
$ws.Range("D51").Value = "1.147"
$ws.Range("E51").Value = "  -6.17%  "
